# Updated cryptos list values (price + volume change) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.182.08"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.18"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.51"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.99"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.811.37"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.608.41"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.47"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.173.32"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.83"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.76"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.95"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.04"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.15"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.411.58"
$ws.Range("E33").Value = "  +8.14%  "
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.88"
$ws.Range("E40").Value = "  +4.84%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.944"
$ws.Range("E42").Value = "  -13.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.763"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.723.06"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.07"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.08"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -2.65%  "
